$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing AntalTest / AntalOmikron figures (rows 8-11) ---
$ws.Range("B8").Value = 5181
$ws.Range("C8").Value = 24

$ws.Range("B9").Value = 4267
$ws.Range("C9").Value = 77

$ws.Range("B10").Value = 4294
$ws.Range("C10").Value = 62

$ws.Range("B11").Value = 4946
$ws.Range("C11").Value = 75

# --- Fill in previously-empty rows 12-14 ---
$ws.Range("B12").Value = 4992
$ws.Range("C12").Value = 110

$ws.Range("B13").Value = 3544
$ws.Range("C13").Value = 110

$ws.Range("B14").Value = 232
$ws.Range("C14").Value = 10

# --- New "Ratio" column D ---
$ws.Range("D1").Value = "Ratio"

$ws.Range("D2").Formula = "=100*C2/B2"
$ws.Range("D3").Formula = "=100*C3/B3"
$ws.Range("D4").Formula = "=100*C4/B4"
$ws.Range("D5").Formula = "=100*C5/B5"
$ws.Range("D6").Formula = "=100*C6/B6"
$ws.Range("D7").Formula = "=100*C7/B7"
$ws.Range("D8").Formula = "=100*C8/B8"
$ws.Range("D9").Formula = "=100*C9/B9"
$ws.Range("D10").Formula = "=100*C10/B10"
$ws.Range("D11").Formula = "=100*C11/B11"
$ws.Range("D12").Formula = "=100*C12/B12"
$ws.Range("D13").Formula = "=100*C13/B13"
$ws.Range("D14").Formula = "=100*C14/B14"

# --- Move the active selection to C16, matching the committed view state ---
$ws.Range("C16").Select()
